# Update the dSF column (F) values to reflect the repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -1
    6  = -3
    7  = -8
    8  = -6
    9  = -11
    10 = -5
    11 = -3
    12 = -6
    13 = 5
    18 = 0
    20 = 2
    23 = 1
    28 = 1
    30 = -2
    31 = -8
    34 = 3
    35 = 0
    40 = 0
    42 = 0
    45 = 0
    46 = 4
    52 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
